$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''68.665.12'
$ws.Range('E2').Value = '''  -1.16%  '
$ws.Range('D3').Value = '''3.863.47'
$ws.Range('E3').Value = '''  -1.94%  '
$ws.Range('D4').Value = '''1.00'
$ws.Range('E4').Value = '''  +0.06%  '
$ws.Range('D5').Value = '''523.54'
$ws.Range('E5').Value = '''  +6.74%  '
$ws.Range('D6').Value = '''141.02'
$ws.Range('E6').Value = '''  -4.16%  '
$ws.Range('D7').Value = '''0.608'
$ws.Range('E7').Value = '''  -2.57%  '
$ws.Range('D8').Value = '''0.999'
$ws.Range('E8').Value = '''  +0.15%  '
$ws.Range('D9').Value = '''0.713'
$ws.Range('E9').Value = '''  -3.16%  '
$ws.Range('E10').Value = '''  -6.23%  '
$ws.Range('D11').Value = '''0.0000321'
$ws.Range('E11').Value = '''  -7.90%  '
$ws.Range('D12').Value = '''41.67'
$ws.Range('E12').Value = '''  -3.30%  '
$ws.Range('D13').Value = '''10.40'
$ws.Range('E13').Value = '''  -0.84%  '
$ws.Range('D14').Value = '''4.476.51'
$ws.Range('E14').Value = '''  -2.03%  '
$ws.Range('B15').Value = 'Chainlink'
$ws.Range('C15').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D15').Value = '''21.27'
$ws.Range('E15').Value = '''  +6.86%  '
$ws.Range('B16').Value = 'WrappedEther'
$ws.Range('C16').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D16').Value = '''3.851.34'
$ws.Range('E16').Value = '''  -2.44%  '
$ws.Range('D17').Value = '''14.10'
$ws.Range('E17').Value = '''  -1.12%  '
$ws.Range('E18').Value = '''  -2.14%  '
$ws.Range('E19').Value = '''  +2.62%  '
$ws.Range('D20').Value = '''68.659.18'
$ws.Range('E20').Value = '''  -1.07%  '
$ws.Range('D21').Value = '''416.36'
$ws.Range('E21').Value = '''  -5.44%  '
$ws.Range('D22').Value = '''3.51'
$ws.Range('E22').Value = '''  +1.53%  '
$ws.Range('D23').Value = '''14.00'
$ws.Range('E23').Value = '''  -3.89%  '
$ws.Range('D24').Value = '''86.85'
$ws.Range('E24').Value = '''  -3.01%  '
$ws.Range('D25').Value = '''4.00'
$ws.Range('E25').Value = '''  +7.18%  '
$ws.Range('D26').Value = '''11.58'
$ws.Range('E26').Value = '''  -4.19%  '
$ws.Range('E27').Value = '''  -5.65%  '
$ws.Range('D28').Value = '''35.63'
$ws.Range('E28').Value = '''  -4.58%  '
$ws.Range('D29').Value = '''13.37'
$ws.Range('E29').Value = '''  -1.04%  '
$ws.Range('D30').Value = '''679.73'
$ws.Range('E30').Value = '''  -3.98%  '
$ws.Range('E31').Value = '''  -5.29%  '
$ws.Range('D32').Value = '''6.83'
$ws.Range('E32').Value = '''  +12.68%  '
$ws.Range('E33').Value = '''  -3.55%  '
$ws.Range('E34').Value = '''  +9.11%  '
$ws.Range('D35').Value = '''0.447'
$ws.Range('E35').Value = '''  -5.72%  '
$ws.Range('D36').Value = '''39.58'
$ws.Range('E36').Value = '''  -3.08%  '
$ws.Range('B37').Value = 'PEPE'
$ws.Range('C37').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D37').Value = '''0.0₃0837'
$ws.Range('E37').Value = '''  -7.55%  '
$ws.Range('B38').Value = 'ThetaToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D38').Value = '''3.51'
$ws.Range('E38').Value = '''  +14.09%  '
$ws.Range('D39').Value = '''0.149'
$ws.Range('E39').Value = '''  -0.79%  '
$ws.Range('E40').Value = '''  +0.14%  '
$ws.Range('E41').Value = '''  -0.11%  '
$ws.Range('D42').Value = '''0.0476'
$ws.Range('D43').Value = '''3.15'
$ws.Range('E43').Value = '''  +4.60%  '
$ws.Range('D44').Value = '''2.80'
$ws.Range('E44').Value = '''  -5.00%  '
$ws.Range('D45').Value = '''3.38'
$ws.Range('E45').Value = '''  +1.02%  '
$ws.Range('E46').Value = '''  -1.39%  '
$ws.Range('D47').Value = '''3.01'
$ws.Range('E47').Value = '''  -2.02%  '
$ws.Range('E48').Value = '''  +12.85%  '
$ws.Range('E49').Value = '''  +0.94%  '
$ws.Range('B50').Value = 'BabyDogeCoin'
$ws.Range('C50').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D50').Value = '''0.0₆0340'
$ws.Range('E50').Value = '''  -7.06%  '
$ws.Range('B51').Value = 'LidoDAOToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D51').Value = '''3.28'
$ws.Range('E51').Value = '''  -2.98%  '
